$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 3
$ws.Range("F3").Value = 257
$ws.Range("F5").Value = 132
$ws.Range("F8").Value = 4813
$ws.Range("F9").Value = 4813
$ws.Range("F10").Value = 18
$ws.Range("F12").Value = 458
$ws.Range("F14").Value = 1102
$ws.Range("F15").Value = 634
$ws.Range("F16").Value = 4420
$ws.Range("F17").Value = 174
$ws.Range("F18").Value = 175
$ws.Range("F21").Value = 3536
$ws.Range("F25").Value = 3220
$ws.Range("F26").Value = 140
$ws.Range("F29").Value = 157
$ws.Range("F30").Value = 201
$ws.Range("F31").Value = 180
$ws.Range("F32").Value = 88
$ws.Range("F33").Value = 67
$ws.Range("F34").Value = 30
$ws.Range("F37").Value = 5620
$ws.Range("F38").Value = 876
$ws.Range("F40").Value = 86
$ws.Range("F41").Value = 959
$ws.Range("F43").Value = 1137
$ws.Range("F44").Value = 510
$ws.Range("F46").Value = 2020
$ws.Range("F47").Value = 302
$ws.Range("F48").Value = 71
$ws.Range("F49").Value = 709
$ws.Range("F50").Value = 863

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 41
$ws.Range("F18").Value = 5
$ws.Range("F24").Value = 749

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 205

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 205
$ws.Range("F3").Value = 257
$ws.Range("F6").Value = 132
$ws.Range("F9").Value = 4813
$ws.Range("F10").Value = 4813
$ws.Range("F11").Value = 18
$ws.Range("F12").Value = 41
$ws.Range("F15").Value = 458
$ws.Range("F16").Value = 1102
$ws.Range("F17").Value = 634
$ws.Range("F18").Value = 4420
$ws.Range("F19").Value = 174
$ws.Range("F20").Value = 175
$ws.Range("F23").Value = 3536
$ws.Range("F24").Value = 3220
$ws.Range("F25").Value = 140
$ws.Range("F27").Value = 157
$ws.Range("F28").Value = 201
$ws.Range("F29").Value = 180
$ws.Range("F30").Value = 88
$ws.Range("F31").Value = 67
$ws.Range("F36").Value = 5620
$ws.Range("F38").Value = 876
$ws.Range("F42").Value = 86
$ws.Range("F43").Value = 959
$ws.Range("F45").Value = 1137
$ws.Range("F46").Value = 510
$ws.Range("F47").Value = 2020
$ws.Range("F48").Value = 302
$ws.Range("F49").Value = 863
